# Add a "Training Required" column (column D) to the "Instal days by Model"
# worksheet's Table1, and populate it with TRUE/FALSE values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instal days by Model")

# Resize the table (Table1) to include the new column D
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D40"))

# Header for the new column
$ws.Range("D1").Value = "Training Required"

# Rows where "Training Required" is FALSE; all other data rows are TRUE
$falseRows = @(16, 23, 38, 39)

for ($r = 2; $r -le 40; $r++) {
    if ($falseRows -contains $r) {
        $ws.Cells.Item($r, 4).Value = $false
    } else {
        $ws.Cells.Item($r, 4).Value = $true
    }
}

# Match the width of the other data columns and select D17 (matching the
# saved view state in the source workbook)
$ws.Columns.Item(4).ColumnWidth = 30.6
$ws.Range("D17").Select()
